$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.048.65'
$ws.Range("E2").Value = '  -1.49%  '

$ws.Range("D3").Value = '2.379.62'
$ws.Range("E3").Value = '  +3.33%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.88%  '

$ws.Range("E7").Value = '  -0.94%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("E9").Value = '  -4.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.34%  '

$ws.Range("E11").Value = '  -2.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.11%  '

$ws.Range("E13").Value = '  -0.25%  '

$ws.Range("D14").Value = '2.738.50'
$ws.Range("E14").Value = '  +3.11%  '

$ws.Range("D15").Value = '2.368.01'
$ws.Range("E15").Value = '  +2.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.825'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.75'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.00%  '

$ws.Range("D18").Value = '45.929.35'
$ws.Range("E18").Value = '  -1.65%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.88%  '

$ws.Range("D20").Value = '0.0₃0953'
$ws.Range("E20").Value = '  +0.60%  '

$ws.Range("E21").Value = '  -0.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("E26").Value = '  -1.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -11.56%  '

$ws.Range("E28").Value = '  -3.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("E30").Value = '  +19.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.02%  '

$ws.Range("E32").Value = '  +6.86%  '

$ws.Range("E33").Value = '  -4.77%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '147.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0775'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.114'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.31%  '

$ws.Range("E37").Value = '  +6.27%  '

$ws.Range("E38").Value = '  -2.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.79%  '

$ws.Range("E41").Value = '  -1.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.43%  '

$ws.Range("D43").Value = '1.931.74'
$ws.Range("E43").Value = '  +3.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.32%  '

$ws.Range("E48").Value = '  -4.91%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '99.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("D50").Value = '2.608.61'
$ws.Range("E50").Value = '  +3.00%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '68.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.15%  '
